$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.766.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.20%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.435.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.88%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'570.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.39%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'146.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.67%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.76%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.111"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.32%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.84%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.16%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'26.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D15").Value = "'2.875.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.74%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'62.635.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.16%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.435.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'11.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.47%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +3.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'324.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.10%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.36%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.02%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +4.76%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'67.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'618.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +8.82%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.54%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +9.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.556.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.89%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +3.59%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.06%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +4.62%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.55%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.12%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +3.17%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.13%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.95%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'5.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.31%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'18.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.17%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.63%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'148.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'2.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +14.19%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.23%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.59%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +2.44%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +1.79%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'20.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.44%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +2.87%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.37%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +4.41%  "
$ws.Range("E51").Style = "Normal"
